# Split the single "Tabelle1" timesheet sheet into two sheets:
#   1. "Logging"   - a small new sheet that stores config data
#                    (carryover row/column pointers)
#   2. "Timesheet" - the original timesheet content (renamed, unchanged)
#
# Final tab order: Logging, Timesheet  (Timesheet = active tab)

$wb = $excel.ActiveWorkbook

# The workbook currently has a single sheet - grab it and rename it to
# "Timesheet"; its data/formatting/drawing stay untouched.
$timesheet = $wb.Worksheets.Item(1)
$timesheet.Name = "Timesheet"

# Insert a brand-new sheet before the Timesheet sheet for the "Logging"
# data and give it the right name / tab order.
$logging = $wb.Worksheets.Add($timesheet)
$logging.Name = "Logging"

# Populate the Logging sheet with the carryover config values.
$logging.Range("B1").Value = "carryover"
$logging.Range("A2").Value = "row"
$logging.Range("B2").Value = 35
$logging.Range("A3").Value = "column"
$logging.Range("B3").Value = 10

# Re-fetch the Timesheet sheet by name - after Add() the old $timesheet
# variable no longer routes writes to the right physical sheet.
$timesheet = $wb.Worksheets.Item("Timesheet")

# Update the Timesheet's selection to match the new authored state.
$timesheet.Range("J35").Select()

# Make "Timesheet" (2nd tab) the active sheet.
$timesheet.Activate()
